$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phone")

# New inventory row: 150 Jio Phones received (S.No. 5), matching the
# existing rows' layout (Qty * Rate = Total Amt Debited).
$ws.Range("A7").Value = 5
$ws.Range("D7").Value = "150 Jio Phone"
$ws.Range("E7").Value = 150
$ws.Range("F7").Value = 1500
$ws.Range("G7").Formula = "=E7*F7"

# Move / extend the active selection as in the authored edit.
[void]$ws.Range("A6:A7").Select()
